$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 608
$ws.Range("K3").Value = 565
$ws.Range("K4").Value = 128
$ws.Range("K5").Value = 34
$ws.Range("K6").Value = 803
$ws.Range("K7").Value = 2138

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K6").Value = 48
$ws.Range("K7").Value = 135

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 25
$ws.Range("K7").Value = 39

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K4").Value = 3
$ws.Range("K6").Value = 32
$ws.Range("K7").Value = 90

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 14
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 42

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 13
$ws.Range("K3").Value = 22
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 64

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 14
$ws.Range("K3").Value = 12
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 18
$ws.Range("K3").Value = 15
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K4").Value = 7
$ws.Range("K8").Value = 135
$ws.Range("K9").Value = 11
$ws.Range("K10").Value = 9
$ws.Range("K11").Value = 53
$ws.Range("K13").Value = 6
$ws.Range("K19").Value = 50
$ws.Range("K20").Value = 58
$ws.Range("K25").Value = 9
$ws.Range("K27").Value = 22
$ws.Range("K29").Value = 113
$ws.Range("K33").Value = 90
$ws.Range("K34").Value = 13
$ws.Range("K37").Value = 64
$ws.Range("K42").Value = 73
$ws.Range("K47").Value = 13
$ws.Range("K51").Value = 32
$ws.Range("K52").Value = 52
$ws.Range("K63").Value = 9
$ws.Range("K65").Value = 58
$ws.Range("K67").Value = 89
$ws.Range("K73").Value = 22
$ws.Range("K79").Value = 53
$ws.Range("K80").Value = 7
$ws.Range("K83").Value = 39
$ws.Range("K84").Value = 14
$ws.Range("K88").Value = 32
$ws.Range("K89").Value = 37
$ws.Range("K90").Value = 18
$ws.Range("K91").Value = 20
$ws.Range("K94").Value = 26
$ws.Range("K95").Value = 42
$ws.Range("K96").Value = 35
$ws.Range("K98").Value = 14
$ws.Range("K99").Value = 49
$ws.Range("K100").Value = 3
$ws.Range("K101").Value = 2138

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 27
$ws.Range("K3").Value = 27
$ws.Range("K4").Value = 7
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 89

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K2").Value = 5
$ws.Range("K7").Value = 14

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 29
$ws.Range("K3").Value = 36
$ws.Range("K7").Value = 113

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 13
$ws.Range("K3").Value = 17
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 24
$ws.Range("K6").Value = 32
$ws.Range("K7").Value = 73

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("I4").Value = 3
$ws.Range("I5").Value = 6

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K3").Value = 2
$ws.Range("K7").Value = 9

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K6").Value = 4
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 53

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("I2").Value = 1
$ws.Range("I7").Value = 3

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K5").Value = 1
$ws.Range("K7").Value = 13

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 9
$ws.Range("K6").Value = 26

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 9

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K2").Value = 4
$ws.Range("K7").Value = 13

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K5").Value = 11
$ws.Range("K6").Value = 14

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 53

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 11

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 8
$ws.Range("K6").Value = 22

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K3").Value = 3
$ws.Range("K4").Value = 1
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K3").Value = 15
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 37

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K2").Value = 4
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 22

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K2").Value = 5
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 18

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K3").Value = 9
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K5").Value = 4
$ws.Range("K6").Value = 7

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 52

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K5").Value = 4
$ws.Range("K6").Value = 7
